# Auto-generated: apply cryptos.xlsx data refresh per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.327.96'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.44%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.911.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.93%  '
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.83'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.661'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +6.16%  '
$ws.Range('E7').Value = '  -0.38%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '41.19'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.31%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.349'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +5.90%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '52.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.74%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0718'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0993'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.190.54'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '12.09'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.31%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.702'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.900.95'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.00%  '
$ws.Range('E17').Value = '  +3.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '35.333.64'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '72.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.45%  '
$ws.Range('D20').Value = '0.0₃0823'
$ws.Range('E20').Value = '  +3.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '239.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.37%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '12.49'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('E23').Value = '  +2.20%  '
$ws.Range('E24').Value = '  -0.39%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.36'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +23.57%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.29'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.76'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.51'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.45'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.62%  '
$ws.Range('E30').Value = '  +2.29%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.16'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.06%  '
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('B33').Value = 'BinanceUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.37%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.936'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.12%  '
$ws.Range('E35').Value = '  +2.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.74'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -4.16%  '
$ws.Range('E37').Value = '  +1.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.34'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.75%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.12'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.54%  '
$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0659'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +11.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0209'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '16.32'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '90.03'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.18%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.338.11'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +3.32%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '47.99'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +38.14%  '
$ws.Range('B47').Value = 'HuobiToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.41'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('B48').Value = 'MXToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.58'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.097.39'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.82%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0703'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.48%  '
